# Version 7 - door sensor tested
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update input values on column H (STM32+RFM95+reed_switch)
# H9  = Consumption Off (ma): 0.002 -> 0.003
$ws.Range("H9").Value = 0.003
# H12 = Battery size (mAh): 300 -> 1400
$ws.Range("H12").Value = 1400

# Recalculate all dependent formulas (H19, H20, H23, H24, etc.)
$excel.CalculateFullRebuild()
$wb.Application.Calculate()

# Update the active cell / selection to H13, matching the saved view state
$ws.Activate()
$ws.Range("H13").Select()
